# "added wireframes and updated moscow"
# Update the MoSCoW requirements list on Sheet1:
#   - Insert a new row for a "sound" requirement under the "should" section
#   - Remove the old "sound" row further down (it effectively moved up)
#   - Rename the "picture based rock paper scissors" idea into
#     "Color recognicion rock paper scissors"
#   - Leave the selection on B22, where the user last clicked

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the old row 12 ("could") so a new "should" item
# (sound) can be added at row 11.
$ws.Rows(11).EntireRow.Insert()

# The old "sound" entry (now shifted to row 15) is no longer needed there
# since it now lives at B11 instead - remove that row entirely.
$ws.Rows(15).EntireRow.Delete()

# Rename the "picture based rock paper scissors" requirement.
$ws.Range("B9").Value = "Color recognicion rock paper scissors"

# Fill in the newly inserted row with the relocated "sound" requirement.
$ws.Range("B11").Value = "sound"

# Restore the last active selection.
$ws.Range("B22").Select()
